$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.071.36'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '2.426.66'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '89.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.74%  '
$ws.Range("E7").Value = '  -2.28%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.46%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0833'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.06%  '
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '31.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.87%  '
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").Value = '2.802.45'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '2.422.19'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.771'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").Value = '41.006.00'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("E19").Value = '  -3.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.30%  '
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.30%  '
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.10%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0746'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.94%  '
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.45%  '
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.60%  '
$ws.Range("E42").Value = '  -7.48%  '
$ws.Range("D43").Value = '1.990.33'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("E45").Value = '  -3.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.92%  '
$ws.Range("D48").Value = '2.665.25'
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '95.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.49%  '
